$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.750.05'
$ws.Range("E2").Value = '  -6.21%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.291.71'
$ws.Range("E3").Value = '  -6.77%  '

$ws.Range("E4").Value = '  +0.25%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '543.40'
$ws.Range("E5").Value = '  -3.12%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '168.62'
$ws.Range("E6").Value = '  -10.31%  '

$ws.Range("E7").Value = '  -3.81%  '

$ws.Range("E8").Value = '  -0.03%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.290.25'
$ws.Range("E9").Value = '  -6.59%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.607'
$ws.Range("E10").Value = '  -4.33%  '

$ws.Range("E11").Value = '  -6.72%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '53.14'
$ws.Range("E12").Value = '  -3.33%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000260'
$ws.Range("E13").Value = '  -5.47%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '8.77'
$ws.Range("E14").Value = '  -6.17%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.834.45'
$ws.Range("E15").Value = '  -6.25%  '

$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.321.43'
$ws.Range("E16").Value = '  -5.99%  '

$ws.Range("B17").Value = 'TRON'
$ws.Range("C17").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.116'
$ws.Range("E17").Value = '  -4.69%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '17.38'
$ws.Range("E18").Value = '  -6.95%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '62.763.78'
$ws.Range("E19").Value = '  -6.19%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.38'
$ws.Range("E20").Value = '  -5.97%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.956'
$ws.Range("E21").Value = '  -4.41%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '397.42'
$ws.Range("E22").Value = '  -5.97%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.96'
$ws.Range("E23").Value = '  -3.25%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.20'
$ws.Range("E24").Value = '  +0.72%  '

$ws.Range("B25").Value = 'Litecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '81.34'
$ws.Range("E25").Value = '  -5.47%  '

$ws.Range("B26").Value = 'InternetComputer(DFINITY)'
$ws.Range("C26").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.92'
$ws.Range("E26").Value = '  +5.23%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.58'
$ws.Range("E27").Value = '  -4.34%  '

$ws.Range("E28").Value = '  -7.58%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.52'
$ws.Range("E29").Value = '  -6.41%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '28.65'
$ws.Range("E30").Value = '  -5.71%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.48'
$ws.Range("E31").Value = '  -2.69%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '577.29'
$ws.Range("E32").Value = '  -9.36%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.13'
$ws.Range("E33").Value = '  -5.22%  '

$ws.Range("E34").Value = '  -6.67%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '57.46'
$ws.Range("E35").Value = '  -4.87%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.147'
$ws.Range("E36").Value = '  -0.47%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.999'
$ws.Range("E37").Value = '  -0.03%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '35.14'
$ws.Range("E38").Value = '  -8.32%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.42'
$ws.Range("E39").Value = '  +1.34%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.360'
$ws.Range("E40").Value = '  -6.94%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0₃0720'
$ws.Range("E41").Value = '  -13.25%  '

$ws.Range("B42").Value = 'FirstDigitalUSD'
$ws.Range("C42").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("E42").Value = '  +0.37%  '

$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.082.11'
$ws.Range("E43").Value = '  -1.39%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.75'
$ws.Range("E44").Value = '  -3.93%  '

$ws.Range("B45").Value = 'ApeXProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.15'
$ws.Range("E45").Value = '  -5.72%  '

$ws.Range("B46").Value = 'Fetch.AI'
$ws.Range("C46").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.43'
$ws.Range("E46").Value = '  -7.94%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0394'
$ws.Range("E47").Value = '  -5.61%  '

$ws.Range("E48").Value = '  -6.71%  '

$ws.Range("E49").Value = '  -4.83%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '130.79'
$ws.Range("E50").Value = '  -6.12%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.88'
$ws.Range("E51").Value = '  -7.62%  '
